$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 233.66667
$ws.Range("I55").Value = 200
$ws.Range("K55").Value = 200
$ws.Range("M55").Value = 14
$ws.Range("H62").Value = 42548.07
$ws.Range("I62").Value = 6791.3335
$ws.Range("K62").Value = 6791.3335
$ws.Range("M62").Value = -6167.3335
$ws.Range("H65").Value = 42548.07
$ws.Range("I65").Value = 6791.3335
$ws.Range("K65").Value = 33956.6675
$ws.Range("M65").Value = -30836.6675
$ws.Range("H80").Value = 1033507.6
$ws.Range("J80").Value = 637.1111
$ws.Range("L80").Value = 1911.3333
$ws.Range("N80").Value = -3907.3333
$ws.Range("H83").Value = 1033507.6
$ws.Range("J83").Value = 637.1111
$ws.Range("L83").Value = 5733.9999
$ws.Range("N83").Value = -15717.9999
$ws.Range("H86").Value = 7734865.5
$ws.Range("I86").Value = 2078.5715
$ws.Range("J86").Value = 10583787
$ws.Range("K86").Value = 2078.5715
$ws.Range("L86").Value = 10583787
$ws.Range("M86").Value = -955.5715
$ws.Range("N86").Value = -10586033
$ws.Range("H89").Value = 7734865.5
$ws.Range("I89").Value = 2078.5715
$ws.Range("J89").Value = 10583787
$ws.Range("K89").Value = 10392.8575
$ws.Range("L89").Value = 52918935
$ws.Range("M89").Value = -4776.8575
$ws.Range("N89").Value = -52930167
$ws.Range("H103").Value = 1110.8889
$ws.Range("I103").Value = 907
$ws.Range("J103").Value = 1169.1428
$ws.Range("K103").Value = 2721
$ws.Range("L103").Value = 3507.4284
$ws.Range("M103").Value = -2135
$ws.Range("N103").Value = -4679.428400000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11120.458
$ws.Range("I32").Value = 6840.136
$ws.Range("J32").Value = 27738.176
$ws.Range("K32").Value = 6840.136
$ws.Range("L32").Value = 27738.176
$ws.Range("M32").Value = -6553.136
$ws.Range("N32").Value = -28312.176
$ws.Range("H45").Value = 618160.9
$ws.Range("I45").Value = 794228.3
$ws.Range("K45").Value = 794228.3
$ws.Range("M45").Value = -793851.3
$ws.Range("H61").Value = 3904.8333
$ws.Range("I61").Value = 2717.6667
$ws.Range("K61").Value = 2717.6667
$ws.Range("M61").Value = -2505.6667
$ws.Range("H97").Value = 1446.7179
$ws.Range("I97").Value = 1130.6471
$ws.Range("K97").Value = 1130.6471
$ws.Range("M97").Value = -634.6470999999999
$ws.Range("H102").Value = 4611160.5
$ws.Range("I102").Value = 5351353.5
$ws.Range("K102").Value = 5351353.5
$ws.Range("M102").Value = -5349731.5
$ws.Range("H136").Value = 3904.8333
$ws.Range("I136").Value = 2717.6667
$ws.Range("K136").Value = 8153.000100000001
$ws.Range("M136").Value = -5603.000100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3456.1667
$ws.Range("I99").Value = 2578.625
$ws.Range("J99").Value = 5211.25
$ws.Range("K99").Value = 2578.625
$ws.Range("L99").Value = 5211.25
$ws.Range("M99").Value = -1080.625
$ws.Range("N99").Value = -8207.25
$ws.Range("H140").Value = 99999.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 99999.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 99999.5
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -110359.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3248.577
$ws.Range("I31").Value = 2572.2632
$ws.Range("K31").Value = 2572.2632
$ws.Range("M31").Value = -2277.2632
$ws.Range("H34").Value = 3248.577
$ws.Range("I34").Value = 2572.2632
$ws.Range("K34").Value = 2572.2632
$ws.Range("M34").Value = -2370.2632
$ws.Range("H58").Value = 2856.697
$ws.Range("I58").Value = 1575
$ws.Range("K58").Value = 1575
$ws.Range("M58").Value = -1372
$ws.Range("H136").Value = 2856.697
$ws.Range("I136").Value = 1575
$ws.Range("K136").Value = 4725
$ws.Range("M136").Value = -2175
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 928925.8
$ws.Range("I122").Value = 1853851.6
$ws.Range("K122").Value = 5561554.800000001
$ws.Range("M122").Value = -5559104.800000001
$ws.Range("H132").Value = 4270.963
$ws.Range("I132").Value = 4012.72
$ws.Range("J132").Value = 7499
$ws.Range("K132").Value = 12038.16
$ws.Range("L132").Value = 22497
$ws.Range("M132").Value = -9508.16
$ws.Range("N132").Value = -27557
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 420
$ws.Range("J22").Value = 420
$ws.Range("L22").Value = 420
$ws.Range("N22").Value = -1010
$ws.Range("H27").Value = 420
$ws.Range("J27").Value = 420
$ws.Range("L27").Value = 420
$ws.Range("N27").Value = -634
$ws.Range("H46").Value = 2801.4167
$ws.Range("I46").Value = 2158.6
$ws.Range("J46").Value = 3260.5715
$ws.Range("K46").Value = 2158.6
$ws.Range("L46").Value = 3260.5715
$ws.Range("M46").Value = -1970.6
$ws.Range("N46").Value = -3636.5715
$ws.Range("H55").Value = 860.6875
$ws.Range("I55").Value = 424.66666
$ws.Range("J55").Value = 1122.3
$ws.Range("K55").Value = 424.66666
$ws.Range("L55").Value = 1122.3
$ws.Range("M55").Value = -251.66666
$ws.Range("N55").Value = -1468.3
$ws.Range("H61").Value = 1083.7142
$ws.Range("I61").Value = 939.8182
$ws.Range("K61").Value = 939.8182
$ws.Range("M61").Value = -737.8182
$ws.Range("H113").Value = 1083.7142
$ws.Range("I113").Value = 939.8182
$ws.Range("K113").Value = 939.8182
$ws.Range("M113").Value = 1230.1818
$ws.Range("H122").Value = 7125.483
$ws.Range("I122").Value = 7086.579
$ws.Range("K122").Value = 21259.737
$ws.Range("M122").Value = -18809.737
$ws.Range("H136").Value = 5799.2456
$ws.Range("I136").Value = 5826.173
$ws.Range("J136").Value = 5519.2
$ws.Range("K136").Value = 17478.519
$ws.Range("L136").Value = 16557.6
$ws.Range("M136").Value = -14928.519
$ws.Range("N136").Value = -21657.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 1375
$ws.Range("I38").Value = 933.6667
$ws.Range("J38").Value = 2699
$ws.Range("K38").Value = 933.6667
$ws.Range("L38").Value = 2699
$ws.Range("M38").Value = -460.6667
$ws.Range("N38").Value = -3645
$ws.Range("H51").Value = 28000
$ws.Range("J51").Value = 30250
$ws.Range("L51").Value = 30250
$ws.Range("N51").Value = -31270
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H81").Value = 9097126
$ws.Range("I81").Value = 4406.0835
$ws.Range("J81").Value = 20008390
$ws.Range("K81").Value = 8812.166999999999
$ws.Range("L81").Value = 40016780
$ws.Range("M81").Value = -7751.166999999999
$ws.Range("N81").Value = -40018902
$ws.Range("H84").Value = 9097126
$ws.Range("I84").Value = 4406.0835
$ws.Range("J84").Value = 20008390
$ws.Range("K84").Value = 44060.835
$ws.Range("L84").Value = 200083900
$ws.Range("M84").Value = -38756.835
$ws.Range("N84").Value = -200094508
$ws.Range("H107").Value = 1424
$ws.Range("I107").Value = 1095.6
$ws.Range("J107").Value = 3066
$ws.Range("K107").Value = 3286.8
$ws.Range("L107").Value = 9198
$ws.Range("M107").Value = -1366.8
$ws.Range("N107").Value = -13038
